$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Eventuele wijzigingen zijn succesvol aangepast." ->
#    "Eventuele wijzigingen zijn succesvol doorgevoerd." split across
#    three runs (matching the target markup), and drop the bCs direct
#    formatting that paragraph carried.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Remove the whole paragraph (text + its own mark); this merges it away
# and leaves the previous paragraph ("Alle code is ...") as the new last
# paragraph.
$lastPara.Range.Delete()

$prevPara = $d.Paragraphs.Last
$insertionPoint = $prevPara.Range.End
$target = $d.Range($insertionPoint, $insertionPoint)

$newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes" ?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr/>
            <w:r>
              <w:rPr/>
              <w:t xml:space="preserve">Eventuele wijzigingen zijn succesvol </w:t>
            </w:r>
            <w:r>
              <w:rPr/>
              <w:t>doorgevoerd</w:t>
            </w:r>
            <w:r>
              <w:rPr/>
              <w:t>.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($newParagraphXml)

# ------------------------------------------------------------------
# 2) Make the (already-portrait) page orientation explicit on the
#    lone section, matching w:orient="portrait" added to pgSz.
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$sec.PageSetup.Orientation = 0
